$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H17").Value = 1081.9
$wsALC.Range("J17").Value = 1107.2632
$wsALC.Range("L17").Value = 3321.7896
$wsALC.Range("N17").Value = -3657.7896
$wsALC.Range("H31").Value = 2032.125
$wsALC.Range("I31").Value = 1608.8572
$wsALC.Range("K31").Value = 4826.571599999999
$wsALC.Range("M31").Value = -4596.571599999999
$wsALC.Range("H40").Value = 2998
$wsALC.Range("I40").Value = 2998
$wsALC.Range("K40").Value = 2998
$wsALC.Range("M40").Value = -2823
$wsALC.Range("H64").Value = 8711.333000000001
$wsALC.Range("I64").Value = 4628.857
$wsALC.Range("K64").Value = 4628.857
$wsALC.Range("M64").Value = -4380.857
$wsALC.Range("H67").Value = 8711.333000000001
$wsALC.Range("I67").Value = 4628.857
$wsALC.Range("K67").Value = 4628.857
$wsALC.Range("M67").Value = -3770.857
$wsALC.Range("H132").Value = 1714.2927
$wsALC.Range("I132").Value = 2151.3215
$wsALC.Range("J132").Value = 773
$wsALC.Range("K132").Value = 6453.9645
$wsALC.Range("L132").Value = 2319
$wsALC.Range("M132").Value = -3923.9645
$wsALC.Range("N132").Value = -7379

# --- ARM ---
$wsARM.Range("H63").Value = 7974.5
$wsARM.Range("I63").Value = 6724.25
$wsARM.Range("J63").Value = 9224.75
$wsARM.Range("K63").Value = 6724.25
$wsARM.Range("L63").Value = 9224.75
$wsARM.Range("M63").Value = -6038.25
$wsARM.Range("N63").Value = -10596.75
$wsARM.Range("H66").Value = 7974.5
$wsARM.Range("I66").Value = 6724.25
$wsARM.Range("J66").Value = 9224.75
$wsARM.Range("K66").Value = 33621.25
$wsARM.Range("L66").Value = 46123.75
$wsARM.Range("M66").Value = -30189.25
$wsARM.Range("N66").Value = -52987.75
$wsARM.Range("H74").Value = 70966.46000000001
$wsARM.Range("I74").Value = 49112
$wsARM.Range("K74").Value = 49112
$wsARM.Range("M74").Value = -48238
$wsARM.Range("H76").Value = 49999
$wsARM.Range("J76").Value = 49999
$wsARM.Range("L76").Value = 49999
$wsARM.Range("N76").Value = -50675
$wsARM.Range("H77").Value = 70966.46000000001
$wsARM.Range("I77").Value = 49112
$wsARM.Range("K77").Value = 245560
$wsARM.Range("M77").Value = -241192
$wsARM.Range("H79").Value = 49999
$wsARM.Range("J79").Value = 49999
$wsARM.Range("L79").Value = 49999
$wsARM.Range("N79").Value = -52339
$wsARM.Range("H80").Value = 54991.25
$wsARM.Range("I80").Value = 0
$wsARM.Range("K80").Value = 0
$wsARM.Range("M80").ClearContents()
$wsARM.Range("H83").Value = 54991.25
$wsARM.Range("I83").Value = 0
$wsARM.Range("K83").Value = 0
$wsARM.Range("M83").ClearContents()
$wsARM.Range("H107").Value = 59998.25
$wsARM.Range("J107").Value = 59998.25
$wsARM.Range("L107").Value = 59998.25
$wsARM.Range("N107").Value = -67678.25

# --- BSM ---
$wsBSM.Range("H35").Value = 64355.9
$wsBSM.Range("J35").Value = 79999
$wsBSM.Range("L35").Value = 79999
$wsBSM.Range("N35").Value = -80619
$wsBSM.Range("H86").Value = 55558988
$wsBSM.Range("I86").Value = 76926270
$wsBSM.Range("K86").Value = 76926270
$wsBSM.Range("M86").Value = -76925147
$wsBSM.Range("H89").Value = 55558988
$wsBSM.Range("I89").Value = 76926270
$wsBSM.Range("K89").Value = 384631350
$wsBSM.Range("M89").Value = -384625734
$wsBSM.Range("H134").Value = 2109.3262
$wsBSM.Range("I134").Value = 1611.1666
$wsBSM.Range("J134").Value = 3902.7
$wsBSM.Range("K134").Value = 4833.4998
$wsBSM.Range("L134").Value = 11708.1
$wsBSM.Range("M134").Value = -2298.4998
$wsBSM.Range("N134").Value = -16778.1

# --- CRP ---
$wsCRP.Range("H41").Value = 31499.666
$wsCRP.Range("J41").Value = 47499.5
$wsCRP.Range("L41").Value = 47499.5
$wsCRP.Range("N41").Value = -48355.5
$wsCRP.Range("H50").Value = 59998
$wsCRP.Range("J50").Value = 59998
$wsCRP.Range("L50").Value = 59998
$wsCRP.Range("N50").Value = -61248
$wsCRP.Range("H51").Value = 46999.1
$wsCRP.Range("J51").Value = 46999.1
$wsCRP.Range("L51").Value = 46999.1
$wsCRP.Range("N51").Value = -48471.1
$wsCRP.Range("H58").Value = 2442.5264
$wsCRP.Range("I58").Value = 1560.6666
$wsCRP.Range("J58").Value = 5749.5
$wsCRP.Range("K58").Value = 1560.6666
$wsCRP.Range("L58").Value = 5749.5
$wsCRP.Range("M58").Value = -1357.6666
$wsCRP.Range("N58").Value = -6155.5
$wsCRP.Range("H60").Value = 41168.285
$wsCRP.Range("J60").Value = 49998.6
$wsCRP.Range("L60").Value = 49998.6
$wsCRP.Range("N60").Value = -51020.6
$wsCRP.Range("H61").Value = 46999.1
$wsCRP.Range("J61").Value = 46999.1
$wsCRP.Range("L61").Value = 46999.1
$wsCRP.Range("N61").Value = -47695.1
$wsCRP.Range("H62").Value = 175450.17
$wsCRP.Range("J62").Value = 11003.667
$wsCRP.Range("L62").Value = 11003.667
$wsCRP.Range("N62").Value = -12251.667
$wsCRP.Range("H65").Value = 175450.17
$wsCRP.Range("J65").Value = 11003.667
$wsCRP.Range("L65").Value = 55018.335
$wsCRP.Range("N65").Value = -61258.335
$wsCRP.Range("H105").Value = 2015.7273
$wsCRP.Range("I105").Value = 2209.125
$wsCRP.Range("K105").Value = 2209.125
$wsCRP.Range("M105").Value = -462.125
$wsCRP.Range("H107").Value = 33365194
$wsCRP.Range("I107").Value = 45496104
$wsCRP.Range("K107").Value = 45496104
$wsCRP.Range("M107").Value = -45494184
$wsCRP.Range("H134").Value = 2471.2114
$wsCRP.Range("I134").Value = 2530.6736
$wsCRP.Range("J134").Value = 1500
$wsCRP.Range("K134").Value = 7592.0208
$wsCRP.Range("L134").Value = 4500
$wsCRP.Range("M134").Value = -5057.0208
$wsCRP.Range("N134").Value = -9570
$wsCRP.Range("H136").Value = 2442.5264
$wsCRP.Range("I136").Value = 1560.6666
$wsCRP.Range("J136").Value = 5749.5
$wsCRP.Range("K136").Value = 4681.9998
$wsCRP.Range("L136").Value = 17248.5
$wsCRP.Range("M136").Value = -2131.9998
$wsCRP.Range("N136").Value = -22348.5

# --- CUL ---
$wsCUL.Range("H69").Value = 15000
$wsCUL.Range("J69").Value = 15000
$wsCUL.Range("L69").Value = 45000
$wsCUL.Range("N69").Value = -46622
$wsCUL.Range("H72").Value = 15000
$wsCUL.Range("J72").Value = 15000
$wsCUL.Range("L72").Value = 135000
$wsCUL.Range("N72").Value = -143112
$wsCUL.Range("H140").Value = 1644.6666
$wsCUL.Range("I140").Value = 1085.5555
$wsCUL.Range("K140").Value = 3256.6665
$wsCUL.Range("M140").Value = 1923.3335

# --- GSM ---
$wsGSM.Range("H26").Value = 0
$wsGSM.Range("I26").Value = 0
$wsGSM.Range("K26").Value = 0
$wsGSM.Range("M26").ClearContents()
$wsGSM.Range("H50").Value = 0
$wsGSM.Range("I50").Value = 0
$wsGSM.Range("K50").Value = 0
$wsGSM.Range("M50").ClearContents()
$wsGSM.Range("H107").Value = 727.2
$wsGSM.Range("I107").Value = 405.75
$wsGSM.Range("K107").Value = 405.75
$wsGSM.Range("M107").Value = 1514.25
$wsGSM.Range("H113").Value = 1530.2858
$wsGSM.Range("I113").Value = 1502.4
$wsGSM.Range("J113").Value = 1600
$wsGSM.Range("K113").Value = 1502.4
$wsGSM.Range("L113").Value = 1600
$wsGSM.Range("M113").Value = 667.5999999999999
$wsGSM.Range("N113").Value = -5940

# --- LTW ---
$wsLTW.Range("H22").Value = 2676.5557
$wsLTW.Range("I22").Value = 577.75
$wsLTW.Range("J22").Value = 3276.2144
$wsLTW.Range("K22").Value = 577.75
$wsLTW.Range("L22").Value = 3276.2144
$wsLTW.Range("M22").Value = -282.75
$wsLTW.Range("N22").Value = -3866.2144
$wsLTW.Range("H27").Value = 2676.5557
$wsLTW.Range("I27").Value = 577.75
$wsLTW.Range("J27").Value = 3276.2144
$wsLTW.Range("K27").Value = 577.75
$wsLTW.Range("L27").Value = 3276.2144
$wsLTW.Range("M27").Value = -470.75
$wsLTW.Range("N27").Value = -3490.2144
$wsLTW.Range("H39").Value = 14784.857
$wsLTW.Range("I39").Value = 14598.8
$wsLTW.Range("K39").Value = 14598.8
$wsLTW.Range("M39").Value = -14138.8
$wsLTW.Range("H46").Value = 3005.9683
$wsLTW.Range("I46").Value = 2200
$wsLTW.Range("K46").Value = 2200
$wsLTW.Range("M46").Value = -2012
$wsLTW.Range("H61").Value = 2973.4546
$wsLTW.Range("I61").Value = 2379.7144
$wsLTW.Range("K61").Value = 2379.7144
$wsLTW.Range("M61").Value = -2177.7144
$wsLTW.Range("H113").Value = 2973.4546
$wsLTW.Range("I113").Value = 2379.7144
$wsLTW.Range("K113").Value = 2379.7144
$wsLTW.Range("M113").Value = -209.7143999999998
$wsLTW.Range("H132").Value = 2927.0256
$wsLTW.Range("I132").Value = 2579.1072
$wsLTW.Range("K132").Value = 7737.321599999999
$wsLTW.Range("M132").Value = -5207.321599999999

# --- WVR ---
$wsWVR.Range("H24").Value = 30000
$wsWVR.Range("J24").Value = 30000
$wsWVR.Range("L24").Value = 30000
$wsWVR.Range("N24").Value = -30460
$wsWVR.Range("H45").Value = 13267.8
$wsWVR.Range("J45").Value = 13736.625
$wsWVR.Range("L45").Value = 13736.625
$wsWVR.Range("N45").Value = -14718.625
$wsWVR.Range("H104").Value = 30658.908
$wsWVR.Range("J104").Value = 30658.908
$wsWVR.Range("L104").Value = 30658.908
$wsWVR.Range("N104").Value = -37646.908
